# Add a new "Serviced by " column (O) to the Card6 sheet, right after the
# existing "Correction" column (N), and tidy up the "Correction " header
# label (drop its trailing space) to match the new column's naming style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card6")

$lastRow = 12

# --- N column header: "Correction " -> "Correction" -------------------
$ws.Range("N1").Value = "Correction"

# --- N column body cells were truly blank before; they now hold "nan" --
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

# --- New O column header: "Serviced by " (note trailing space) --------
$ws.Range("O1").Value = "Serviced by "
# match the style already used by the other header cells (bold, border,
# centered / top aligned) by copying N1's formatting onto O1
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- New O column body cells: created blank, same as N used to be -----
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Style = "Normal"
}
